$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 8.660232485948974, 0, 14.90378790461981)
    3 = @(1.459612070389937, 10.29869402782916, 0.1575252929769615, 8.660232485948974, 0, 20.57606387714503)
    4 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.997429241610044)
    5 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 7.524616544037286)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
